# Updates cryptos list cell values per the Wed Nov 13 09:31:31 UTC 2024 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) cells are stored as plain text in the workbook (not numbers),
# e.g. "87.567.32" (European-style thousands separators) would otherwise be
# auto-coerced to a number by Excel. Prefixing the value with a literal leading
# apostrophe forces Excel to keep/store it as text, matching the source data.

# Row 2
$ws.Range("D2").Value = "'87.567.32"
$ws.Range("E2").Value = "  -1.20%  "

# Row 3
$ws.Range("D3").Value = "'3.163.82"
$ws.Range("E3").Value = "  -7.11%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").Value = "'207.63"
$ws.Range("E5").Value = "  -5.67%  "

# Row 6
$ws.Range("D6").Value = "'608.81"
$ws.Range("E6").Value = "  -6.65%  "

# Row 7
$ws.Range("D7").Value = "'0.376"
$ws.Range("E7").Value = "  -11.50%  "

# Row 8
$ws.Range("D8").Value = "'0.665"
$ws.Range("E8").Value = "  -0.77%  "

# Row 9
$ws.Range("E9").Value = "  -0.02%  "

# Row 10
$ws.Range("D10").Value = "'3.160.92"
$ws.Range("E10").Value = "  -7.09%  "

# Row 11
$ws.Range("D11").Value = "'0.533"
$ws.Range("E11").Value = "  -14.38%  "

# Row 12
$ws.Range("E12").Value = "  +4.18%  "

# Row 13
$ws.Range("D13").Value = "'0.0000239"
$ws.Range("E13").Value = "  -17.21%  "

# Row 14
$ws.Range("D14").Value = "'3.745.62"
$ws.Range("E14").Value = "  -7.00%  "

# Row 15
$ws.Range("D15").Value = "'5.23"
$ws.Range("E15").Value = "  -6.68%  "

# Row 16
$ws.Range("D16").Value = "'87.305.39"
$ws.Range("E16").Value = "  -1.30%  "

# Row 17
$ws.Range("D17").Value = "'32.10"
$ws.Range("E17").Value = "  -13.05%  "

# Row 18
$ws.Range("D18").Value = "'3.154.12"
$ws.Range("E18").Value = "  -7.45%  "

# Row 19
$ws.Range("E19").Value = "  -0.21%  "

# Row 20
$ws.Range("D20").Value = "'13.40"
$ws.Range("E20").Value = "  -11.19%  "

# Row 21
$ws.Range("D21").Value = "'414.15"
$ws.Range("E21").Value = "  -10.39%  "

# Row 22
$ws.Range("D22").Value = "'8.42"
$ws.Range("E22").Value = "  -13.77%  "

# Row 23
$ws.Range("D23").Value = "'5.05"
$ws.Range("E23").Value = "  -11.09%  "

# Row 24
$ws.Range("D24").Value = "'5.16"
$ws.Range("E24").Value = "  -7.64%  "

# Row 25
$ws.Range("D25").Value = "'11.85"
$ws.Range("E25").Value = "  -7.71%  "

# Row 26
$ws.Range("D26").Value = "'3.333.15"
$ws.Range("E26").Value = "  -6.89%  "

# Row 27
$ws.Range("D27").Value = "'73.22"
$ws.Range("E27").Value = "  -10.21%  "

# Row 28
$ws.Range("E28").Value = "  -11.25%  "

# Row 29
$ws.Range("E29").Value = "  -0.02%  "

# Row 30
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  -0.11%  "

# Row 31
$ws.Range("D31").Value = "'0.156"
$ws.Range("E31").Value = "  -17.04%  "

# Row 32
$ws.Range("D32").Value = "'543.32"
$ws.Range("E32").Value = "  -6.92%  "

# Row 33
$ws.Range("E33").Value = "  -13.94%  "

# Row 34
$ws.Range("E34").Value = "  -17.72%  "

# Row 35
$ws.Range("B35").Value = "PancakeSwap"
$ws.Range("C35").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D35").Value = "'1.84"
$ws.Range("E35").Value = "  -12.94%  "

# Row 36
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D36").Value = "'6.69"
$ws.Range("E36").Value = "  -11.08%  "

# Row 37
$ws.Range("E37").Value = "  -8.35%  "

# Row 38
$ws.Range("D38").Value = "'21.73"
$ws.Range("E38").Value = "  -9.89%  "

# Row 39
$ws.Range("D39").Value = "'21.82"
$ws.Range("E39").Value = "  -0.11%  "

# Row 40
$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  +0.21%  "

# Row 41
$ws.Range("D41").Value = "'2.99"
$ws.Range("E41").Value = "  -5.13%  "

# Row 42
$ws.Range("E42").Value = "  +0.02%  "

# Row 43
$ws.Range("E43").Value = "  -11.06%  "

# Row 44
$ws.Range("D44").Value = "'0.367"
$ws.Range("E44").Value = "  -16.61%  "

# Row 45
$ws.Range("D45").Value = "'148.47"
$ws.Range("E45").Value = "  -6.38%  "

# Row 46
$ws.Range("E46").Value = "  -9.73%  "

# Row 47
$ws.Range("D47").Value = "'43.14"
$ws.Range("E47").Value = "  -7.77%  "

# Row 48
$ws.Range("D48").Value = "'0.124"
$ws.Range("E48").Value = "  -0.61%  "

# Row 49
$ws.Range("E49").Value = "  -15.20%  "

# Row 50
$ws.Range("E50").Value = "  -13.15%  "

# Row 51
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "'23.69"
$ws.Range("E51").Value = "  -9.27%  "
